# Case study 2 (Wang 2024) - feature importance sheet refresh.
# The booster was re-run and the "gain" ranking re-sorted in descending
# order; this rewrites the (feature, gain) pairs in A2:B121 to match the
# newly exported ranking while leaving the header row untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rankedFeatures = @(
    @(2, 'YB/LU', 1572),
    @(3, 'TM/YB', 1462),
    @(4, 'TH/U', 1408),
    @(5, 'TB/DY', 1360),
    @(6, 'CE**EU*EU', 1360),
    @(7, 'ER/TM', 1348),
    @(8, 'DY/HO', 1298),
    @(9, 'Y/LA/TB', 1282),
    @(10, 'EU/GD', 1248),
    @(11, 'Y/CE*/TH', 1156),
    @(12, 'TA/U', 1148),
    @(13, 'HF-U', 1120),
    @(14, 'LU/TH', 1098),
    @(15, 'HO-LU', 1090),
    @(16, 'EU**EU**TA', 1072),
    @(17, 'YB/LU/HF', 1052),
    @(18, 'CE/GD', 1050),
    @(19, 'ER/YB', 1040),
    @(20, 'Y/YB/TH', 1006),
    @(21, 'NB/LA/CE*', 1000),
    @(22, 'CE/HF/TH', 995.9999999999997),
    @(23, 'TI*Y*EU*', 985.9999999999998),
    @(24, 'CE/TA/U', 983.9999999999998),
    @(25, 'PR/EU/TA', 981.9999999999998),
    @(26, 'CE/U', 976.0000000000002),
    @(27, 'PR/ND/EU*', 957.9999999999998),
    @(28, 'CE-ND-TA', 945.9999999999994),
    @(29, 'Y/CE*/EU', 942.0000000000003),
    @(30, 'Y/CE/EU*', 939.9999999999997),
    @(31, 'Y/DY/HF', 936.0000000000001),
    @(32, 'CE/DY/U', 929.9999999999995),
    @(33, 'TI/HF/HF', 926.0000000000002),
    @(34, 'NB/LU/TA', 925.9999999999997),
    @(35, 'TH-U', 920),
    @(36, 'CE/EU*/TH', 910),
    @(37, 'TI*EU*HF', 899.9999999999999),
    @(38, 'NB/CE/CE', 897.9999999999998),
    @(39, 'TI*TI*U', 894.0000000000005),
    @(40, 'EU**HF*TA', 892.0000000000002),
    @(41, 'TI/NB/EU*', 885.9999999999999),
    @(42, 'YB/LU/TA', 881.9999999999995),
    @(43, 'PR/ND/TA', 878.0000000000001),
    @(44, 'CE/EU*/EU*', 866),
    @(45, 'TI/EU*/HF', 854),
    @(46, 'Y/NB', 852.0000000000001),
    @(47, 'TI*TI*CE', 840.0000000000005),
    @(48, 'ER/TM/HF', 837.9999999999999),
    @(49, 'TI*NB*TA', 834),
    @(50, 'TI*TI*TA', 824),
    @(51, 'TI*Y*TA', 818.0000000000003),
    @(52, 'Y/NB/ER', 817.9999999999999),
    @(53, 'HF-TH-U', 816.0000000000001),
    @(54, 'TI*HF*TA', 814.0000000000003),
    @(55, 'GD/TB', 807.9999999999999),
    @(56, 'CE/EU*', 806),
    @(57, 'EU**HF*HF', 805.9999999999998),
    @(58, 'Y/CE/U', 800),
    @(59, 'Y*CE*EU*', 798.0000000000005),
    @(60, 'CE/EU/TH', 792.0000000000001),
    @(61, 'DY-LU-LU', 787.9999999999995),
    @(62, 'ER-LU-LU', 784.0000000000003),
    @(63, 'CE/DY/DY', 784),
    @(64, 'SM/EU/TA', 783.9999999999999),
    @(65, 'ER/YB/HF', 780),
    @(66, 'TI*EU**U', 778),
    @(67, 'SM/GD', 756.0000000000001),
    @(68, 'LA+EU*+TA', 740.0000000000001),
    @(69, 'CE*EU**HF', 736.0000000000002),
    @(70, 'ND-EU-EU*', 731.9999999999998),
    @(71, 'CE/GD/HF', 731.9999999999995),
    @(72, 'CE**ND*HF', 727.9999999999999),
    @(73, 'ER/LU', 723.9999999999999),
    @(74, 'DY/LU', 722.0000000000001),
    @(75, 'CE/TH', 718.0000000000001),
    @(76, 'NB*EU**U', 715.9999999999999),
    @(77, 'NB/HF/U', 708.0000000000002),
    @(78, 'EU/GD/HF', 707.9999999999997),
    @(79, 'EU*/TA/TA', 701.9999999999998),
    @(80, 'Y/HO', 699.9999999999995),
    @(81, 'Y/EU*/U', 698.0000000000002),
    @(82, 'TM/YB/HF', 695.9999999999995),
    @(83, 'HO/TM', 691.9999999999995),
    @(84, 'Y/TM', 688),
    @(85, 'CE/EU*/U', 686),
    @(86, 'NB/HF/TH', 685.9999999999999),
    @(87, 'CE/TA/TH', 676.0000000000002),
    @(88, 'TI/EU*/TA', 676.0000000000001),
    @(89, 'Y/ER', 670),
    @(90, 'HO-LU-TA', 662),
    @(91, 'TI*HF*HF', 660),
    @(92, 'EU/EU*/GD', 656.0000000000001),
    @(93, 'Y*HF*TA', 650.0000000000002),
    @(94, 'HO/ER', 648.0000000000001),
    @(95, 'Y/NB/YB', 634.0000000000001),
    @(96, 'TI/HF', 632.0000000000002),
    @(97, 'Y/SM', 630),
    @(98, 'CE/HF/HF', 628.0000000000002),
    @(99, 'NB-TA-TA', 624.0000000000001),
    @(100, 'SM/GD/HF', 622.0000000000002),
    @(101, 'CE/SM/HF', 621.9999999999999),
    @(102, 'Y/NB/EU*', 611.9999999999999),
    @(103, 'CE/EU/EU*', 610.0000000000001),
    @(104, 'Y/YB', 603.9999999999998),
    @(105, 'TM/LU', 601.9999999999999),
    @(106, 'NB/TH/U', 599.9999999999999),
    @(107, 'Y/NB/CE', 597.9999999999999),
    @(108, 'CE/EU/TA', 586),
    @(109, 'DY-LU-TA', 584.0000000000001),
    @(110, 'TI*TI*DY', 580.0000000000002),
    @(111, 'Y*EU**HF', 572),
    @(112, 'Y/DY', 571.9999999999998),
    @(113, 'Y/TH', 566.0000000000003),
    @(114, 'ND-EU-EU', 564),
    @(115, 'CE/EU/HF', 540.0000000000002),
    @(116, 'TB/DY/HF', 534.0000000000001),
    @(117, 'Y/DY/TA', 531.9999999999998),
    @(118, 'TB/LU/HF', 432.0000000000001),
    @(119, 'NB/HF/TA', 425.9999999999999),
    @(120, 'TA*U', 414.0000000000001),
    @(121, 'NB/TA', 334)
)

foreach ($row in $rankedFeatures) {
    $rowIndex = $row[0]
    $featureName = $row[1]
    $gainValue = $row[2]

    $ws.Cells.Item($rowIndex, 1).Value = $featureName
    $ws.Cells.Item($rowIndex, 2).Value = $gainValue
}
